# Update the division worksheet's answer key.
# The document has a single 5-column table; data rows (with content) sit
# at table rows 1, 5, 9, 13, 17 (rows 2-4, 6-8, etc. are blank "work" rows).
# Replace each cell's text directly via Table.Cell(row, col).Range.Text so
# the existing run formatting (rFonts/sz) is preserved untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "40÷7=5, 5"
$t.Cell(1, 2).Range.Text  = "75÷7=10, 5"
$t.Cell(1, 3).Range.Text  = "82÷3=27, 1"
$t.Cell(1, 4).Range.Text  = "43÷2=21, 1"
$t.Cell(1, 5).Range.Text  = "69÷7=9, 6"

$t.Cell(5, 1).Range.Text  = "68÷4=17, 0"
$t.Cell(5, 2).Range.Text  = "89÷5=17, 4"
$t.Cell(5, 3).Range.Text  = "69÷5=13, 4"
$t.Cell(5, 4).Range.Text  = "21÷9=2, 3"
$t.Cell(5, 5).Range.Text  = "55÷2=27, 1"

$t.Cell(9, 1).Range.Text  = "60÷7=8, 4"
$t.Cell(9, 2).Range.Text  = "14÷7=2, 0"
$t.Cell(9, 3).Range.Text  = "19÷2=9, 1"
$t.Cell(9, 4).Range.Text  = "39÷6=6, 3"
$t.Cell(9, 5).Range.Text  = "34÷9=3, 7"

$t.Cell(13, 1).Range.Text = "11÷8=1, 3"
$t.Cell(13, 2).Range.Text = "64÷2=32, 0"
$t.Cell(13, 3).Range.Text = "26÷3=8, 2"
$t.Cell(13, 4).Range.Text = "34÷5=6, 4"
$t.Cell(13, 5).Range.Text = "31÷8=3, 7"

$t.Cell(17, 1).Range.Text = "51÷2=25, 1"
$t.Cell(17, 2).Range.Text = "29÷2=14, 1"
$t.Cell(17, 3).Range.Text = "43÷8=5, 3"
$t.Cell(17, 4).Range.Text = "90÷7=12, 6"
$t.Cell(17, 5).Range.Text = "56÷6=9, 2"
